$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 113, pushing existing rows 113-153 down to 115-155.
$ws.Rows("113:114").Insert()

# New row 113 data
$ws.Cells.Item(113, 1).Value = 6
$ws.Cells.Item(113, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(113, 3).Value = "Metropolitana"
$ws.Cells.Item(113, 4).Value = 45029
$ws.Cells.Item(113, 5).Value = 13
$ws.Cells.Item(113, 6).Value = "Fruta"
$ws.Cells.Item(113, 7).Value = 100104
$ws.Cells.Item(113, 8).Value = "Frutos de pepita"
$ws.Cells.Item(113, 9).Value = 100104003
$ws.Cells.Item(113, 10).Value = "Membrillo"
$ws.Cells.Item(113, 11).Value = "Champion"
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 6
$ws.Cells.Item(113, 14).Value = 230000
$ws.Cells.Item(113, 15).Value = 230000
$ws.Cells.Item(113, 16).Value = 230000
$ws.Cells.Item(113, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(113, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(113, 19).Value = 511
$ws.Cells.Item(113, 20).Value = 450

# New row 114 data
$ws.Cells.Item(114, 1).Value = 6
$ws.Cells.Item(114, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(114, 3).Value = "Metropolitana"
$ws.Cells.Item(114, 4).Value = 45029
$ws.Cells.Item(114, 5).Value = 13
$ws.Cells.Item(114, 6).Value = "Fruta"
$ws.Cells.Item(114, 7).Value = 100104
$ws.Cells.Item(114, 8).Value = "Frutos de pepita"
$ws.Cells.Item(114, 9).Value = 100104003
$ws.Cells.Item(114, 10).Value = "Membrillo"
$ws.Cells.Item(114, 11).Value = "Champion"
$ws.Cells.Item(114, 12).Value = "Segunda"
$ws.Cells.Item(114, 13).Value = 8
$ws.Cells.Item(114, 14).Value = 180000
$ws.Cells.Item(114, 15).Value = 180000
$ws.Cells.Item(114, 16).Value = 180000
$ws.Cells.Item(114, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(114, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(114, 19).Value = 400
$ws.Cells.Item(114, 20).Value = 450
